$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.119.35"
$ws.Range("E2").Value = "  -6.13%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.554.11"
$ws.Range("E3").Value = "  -2.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "299.14"
$ws.Range("E5").Value = "  -3.13%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.99"
$ws.Range("E6").Value = "  -5.99%  "

$ws.Range("E7").Value = "  -3.20%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("E9").Value = "  -4.80%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.18"
$ws.Range("E10").Value = "  -6.66%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0811"
$ws.Range("E11").Value = "  -3.74%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.76"
$ws.Range("E12").Value = "  -3.96%  "

$ws.Range("E13").Value = "  +1.79%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.949.31"
$ws.Range("E14").Value = "  -2.02%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.552.43"
$ws.Range("E15").Value = "  -2.27%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.875"
$ws.Range("E16").Value = "  -4.01%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.20"
$ws.Range("E17").Value = "  -4.02%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.183.35"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.02"
$ws.Range("E19").Value = "  +2.65%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0982"
$ws.Range("E20").Value = "  -3.09%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.63"
$ws.Range("E21").Value = "  -1.67%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.39"
$ws.Range("E22").Value = "  -0.31%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "260.54"
$ws.Range("E23").Value = "  -10.49%  "

$ws.Range("E24").Value = "  -3.40%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "29.69"
$ws.Range("E25").Value = "  +0.22%  "

$ws.Range("E26").Value = "  -4.71%  "

$ws.Range("E27").Value = "  -0.06%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.04"
$ws.Range("E28").Value = "  -6.63%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.19"
$ws.Range("E29").Value = "  -3.92%  "

$ws.Range("E30").Value = "  -3.68%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.00"
$ws.Range("E31").Value = "  -4.05%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "154.51"
$ws.Range("E32").Value = "  -3.51%  "

$ws.Range("E33").Value = "  -2.70%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.76"
$ws.Range("E34").Value = "  -2.01%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.38"
$ws.Range("E35").Value = "  -6.71%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0798"
$ws.Range("E36").Value = "  -4.93%  "

$ws.Range("E37").Value = "  -4.66%  "

$ws.Range("E38").Value = "  -2.63%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.76"
$ws.Range("E39").Value = "  +6.91%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.34"
$ws.Range("E40").Value = "  +9.73%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.49"
$ws.Range("E41").Value = "  -1.57%  "

$ws.Range("E42").Value = "  -4.70%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.89"
$ws.Range("E43").Value = "  -2.89%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.067.30"
$ws.Range("E44").Value = "  -2.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.998"
$ws.Range("E45").Value = "  -0.08%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "85.86"
$ws.Range("E46").Value = "  -10.58%  "

$ws.Range("E47").Value = "  +3.15%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.804.41"
$ws.Range("E48").Value = "  -2.33%  "

$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.71"
$ws.Range("E49").Value = "  -0.76%  "

$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.76"
$ws.Range("E50").Value = "  -6.86%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "104.36"
$ws.Range("E51").Value = "  -4.62%  "
